$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 21; existing rows 21-88 shift down to 22-89.
$ws.Rows("21:21").Insert()

# Populate the new row (ID 573, apartment name 청계벽산).
$ws.Range("A21").Value = 573
$ws.Range("B21").Value = "청계벽산"

# New cell uses a distinct font (맑은 고딕, 9pt) while keeping the sheet's
# standard vertical-center alignment.
$ws.Range("B21").Font.Name = "맑은 고딕"
$ws.Range("B21").Font.Size = 9

# Row grows slightly taller to fit the new font.
$ws.Rows("21:21").RowHeight = 13.2

# Update the view: selection moves to M23 and the scroll position resets.
[void]$ws.Range("M23").Select()

Write-Host "Applied edit"
